# Update the "Periodo Mora" column (E16:E22) so that the periods are listed
# newest-to-oldest instead of oldest-to-newest, and keep the "Valor Mora"
# (F column) value tied to the correct period: the most recent period (2306)
# now occupies the first data row and keeps its 29333 value, while the
# oldest period (2212) moves to the last data row and keeps its 40000 value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @("2306", "2305", "2304", "2303", "2302", "2301", "2212")
$valores = @(29333, 40000, 40000, 40000, 40000, 40000, 40000)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]   # Column E = Periodo Mora
    $ws.Cells.Item($row, 6).Value = $valores[$i]    # Column F = Valor Mora
}
